$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected, so locked cells (the whole data area) can't be
# written to directly. Temporarily unprotect, make the edits, then
# re-protect so the workbook ends up protected again.
$ws.Unprotect()

# Update the "as of" date in the disclosure note (A9).
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-10 for illustrative purposes only and are subject to change."

# Re-fit row 9's height (writing the multi-line text otherwise stamps an
# explicit custom row height that wasn't present in the original file).
$ws.Rows.Item(9).AutoFit()

# Update the Weight / Percent Change figures.
$ws.Range("D2").Value = 0.255140977909357
$ws.Range("E2").Value = -0.004438374871969963

$ws.Range("D3").Value = 0.2531722655462616
$ws.Range("E3").Value = -0.0007940709370035659

$ws.Range("D4").Value = 0.2452402174721214
$ws.Range("E4").Value = 0.0009492168960609781

$ws.Range("D5").Value = 0.2464465390722599
$ws.Range("E5").Value = -0.01916198262646907

$ws.Range("E6").Value = -0.005823066185326464

# Restore sheet protection.
$ws.Protect()
